# Edit: "Fruta / hortaliza, semanal"
#
# The sheet holds a weekly price table for "Espárragos" at Mercado Mayorista
# Lo Valledor de Santiago. This update inserts 4 new observation rows into
# the existing table (rows 112-137), pushing the subsequent rows down, so
# that the sheet ends up with rows 1:R141 instead of 1:R137.
#
#   - 1 new row is inserted at row 112 (week of 2021-09-28 / serial 44467).
#   - 3 new rows are inserted at rows 129-131 (week of 2022-10-07 / serial 44841).
#
# All rows below each insertion point shift down accordingly; their values
# are left untouched (Excel carries them down automatically).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DataRow {
    # NOTE: use positional parameters here - named parameter binding
    # (-Row / -Values) does not reliably bind in this COM-interop host.
    param(
        [int]$Row,
        [object[]]$Values
    )
    for ($i = 0; $i -lt $Values.Length; $i++) {
        $ws.Cells.Item($Row, $i + 1).Value = $Values[$i]
    }
}

# --- Insert 1 row at 112 ---------------------------------------------------
$ws.Rows.Item(112).Insert()

$row112 = @(
    6,
    "Mercado Mayorista Lo Valledor de Santiago",
    "Metropolitana",
    44467,
    13,
    300000000,
    "Espárragos",
    "Sin especificar",
    "Segunda",
    50,
    1300,
    1300,
    1300,
    '$/caja 10 kilos',
    "Provincia de Linares",
    130,
    10,
    "Hortaliza"
)
Set-DataRow 112 $row112

# --- Insert 3 rows at 129:131 ----------------------------------------------
$ws.Range("129:131").Insert()

$row129 = @(
    6,
    "Mercado Mayorista Lo Valledor de Santiago",
    "Metropolitana",
    44841,
    13,
    300000000,
    "Espárragos",
    "Sin especificar",
    "Banquete",
    670,
    1900,
    2000,
    1952,
    '$/kilo',
    "Provincia de Linares",
    1952,
    1,
    "Hortaliza"
)
Set-DataRow 129 $row129

$row130 = @(
    6,
    "Mercado Mayorista Lo Valledor de Santiago",
    "Metropolitana",
    44841,
    13,
    300000000,
    "Espárragos",
    "Sin especificar",
    "Primera",
    510,
    1600,
    1700,
    1649,
    '$/kilo',
    "Provincia de Linares",
    1649,
    1,
    "Hortaliza"
)
Set-DataRow 130 $row130

$row131 = @(
    6,
    "Mercado Mayorista Lo Valledor de Santiago",
    "Metropolitana",
    44841,
    13,
    300000000,
    "Espárragos",
    "Sin especificar",
    "Segunda",
    370,
    1400,
    1500,
    1454,
    '$/kilo',
    "Provincia de Linares",
    1454,
    1,
    "Hortaliza"
)
Set-DataRow 131 $row131

Write-Host ("Inserted 4 rows; dimension is now " + $ws.UsedRange.Address())
